# Weekly data update: insert two new price records (row 251 and 252) for
# "Ají" at "Feria Lagunitas de Puerto Montt", pushing the existing rows
# 251-320 down to 253-322 (the sheet grows from A1:R320 to A1:R322).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 251; Excel shifts rows 251:320 down to 253:322
$ws.Rows("251:252").Insert()

# --- New row 251 ---------------------------------------------------------
$ws.Cells.Item(251, 1).Value = 4
$ws.Cells.Item(251, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(251, 3).Value = "Los Lagos"
$ws.Cells.Item(251, 4).Value = 44841
$ws.Cells.Item(251, 5).Value = 10
$ws.Cells.Item(251, 6).Value = 100112021
$ws.Cells.Item(251, 7).Value = "Ají"
$ws.Cells.Item(251, 8).Value = "Inferno"
$ws.Cells.Item(251, 9).Value = "Primera"
$ws.Cells.Item(251, 10).Value = 90
$ws.Cells.Item(251, 11).Value = 30000
$ws.Cells.Item(251, 12).Value = 30000
$ws.Cells.Item(251, 13).Value = 30000
$ws.Cells.Item(251, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(251, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(251, 16).Value = 3000
$ws.Cells.Item(251, 17).Value = 10
$ws.Cells.Item(251, 18).Value = "Hortaliza"

# --- New row 252 ---------------------------------------------------------
$ws.Cells.Item(252, 1).Value = 4
$ws.Cells.Item(252, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(252, 3).Value = "Los Lagos"
$ws.Cells.Item(252, 4).Value = 44841
$ws.Cells.Item(252, 5).Value = 10
$ws.Cells.Item(252, 6).Value = 100112021
$ws.Cells.Item(252, 7).Value = "Ají"
$ws.Cells.Item(252, 8).Value = "Inferno"
$ws.Cells.Item(252, 9).Value = "Segunda"
$ws.Cells.Item(252, 10).Value = 90
$ws.Cells.Item(252, 11).Value = 22000
$ws.Cells.Item(252, 12).Value = 22000
$ws.Cells.Item(252, 13).Value = 22000
$ws.Cells.Item(252, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(252, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(252, 16).Value = 2200
$ws.Cells.Item(252, 17).Value = 10
$ws.Cells.Item(252, 18).Value = "Hortaliza"

# Ensure the date cells keep the workbook's date number format (style index 2
# on D2 in the original sheet), matching the other rows in column D.
$ws.Cells.Item(251, 4).NumberFormat = $ws.Cells.Item(253, 4).NumberFormat()
$ws.Cells.Item(252, 4).NumberFormat = $ws.Cells.Item(253, 4).NumberFormat()
